# Update column F (dSF) values on the active sheet.
# This reflects a repull of data for the "dSF" column, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = 3
    4  = 7
    5  = -2
    6  = -5
    7  = 2
    9  = 2
    10 = -3
    11 = 3
    12 = 5
    13 = 1
    14 = -3
    15 = -3
    16 = 2
    17 = 4
    18 = 3
    19 = -3
    20 = -2
    21 = 3
    22 = -2
    23 = 1
    24 = 0
    25 = 1
    26 = -1
    27 = -3
    28 = 1
    29 = 2
    30 = -6
    31 = 1
    32 = -1
    33 = -3
    34 = 1
    35 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
